# Restore C10 value on the active sheet from 18 to 1 (numeric), matching the
# authored commit's change of xl/worksheets/sheet1.xml (<c r="C10"> 18 -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
